# Rotate the comma-separated "Recorded By" values in column G so that the
# last author in the list moves to the front (rotate-right by one element).
# Cells with only a single value are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ",\s*"
    if ($parts.Count -le 1) { continue }

    $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
    $newVal = [string]::Join(", ", $rotated)

    $cell.Value = $newVal
}
